# Update the header in B1 (merged B1:D1) on Sheet1 of the TIMES_Energy_SUN
# workbook: the column order in the label changed from
#   "Unit - activity1  /  Period  /  Process Description1"
# to
#   "Period  /  Unit - activity1  /  Process Description1"
# as part of refreshing the EV profiles source data (RAMP-mobility results).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Leading apostrophe keeps the cell's existing "quote prefix" text style
# (xf index 7 in styles.xml) instead of Excel minting a brand-new style
# record that differs only by the missing quotePrefix flag.
$ws.Range("B1").Value = "'Period  /  Unit - activity1  /  Process Description1"

# The saved view had scrolled/selected deep into the sheet (topLeftCell
# A16, selection A25); reset the view back to the top-left corner, as in
# the saved file (no scroll offset, default selection).
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
